# Remediation-Notice.docx edit
#
# Replaces the first two blank paragraphs that follow the
# "List the remedial actions necessary to afford compliance with the
# relevant instrument:" paragraph with a Jinja/docxtpl loop:
#
#   {% for ra in remediation_actions %}{{ ra.action }} ( Due date: {{ ra.due_date }} )
#   {% endfor %}
#
# The first new paragraph carries bold / 10pt (sz=20) paragraph-mark
# formatting, and its second run (the "{{ ra.action }} ..." text) is
# bold / 10pt as well. The second new paragraph ("{% endfor %}") keeps
# the original, unformatted paragraph mark.

$d = $word.ActiveDocument

# --- Locate the anchor paragraph ("List the remedial actions ...") ---
$anchorIndex = -1
$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*List the remedial actions necessary to afford compliance*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the 'List the remedial actions...' paragraph"
}

# The two paragraphs immediately following the anchor are the blank
# paragraphs that need to become the for-loop / endfor paragraphs.
$forParaIndex = $anchorIndex + 1
$endForParaIndex = $anchorIndex + 2

# --- Paragraph 1: "{% for ra in remediation_actions %}" + bold action text ---
$forPara = $d.Paragraphs.Item($forParaIndex)
$forRange = $forPara.Range

$forXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b w:val="1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">{% for ra in remediation_actions %}</w:t></w:r><w:r><w:rPr><w:b w:val="1"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">{{ ra.action }} ( Due date: {{ ra.due_date }} )</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$forRange.InsertXML($forXml)

# --- Paragraph 2: "{% endfor %}" (keeps its existing plain formatting) ---
$endForPara = $d.Paragraphs.Item($endForParaIndex)
$endForRange = $endForPara.Range
[void]$endForRange.InsertBefore("{% endfor %}")
